# Swap the data held in rows 20 and 21 (two species-observation records
# exchange places), per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row1 = 20
$row2 = 21

# Columns that carry differing data between the two rows.
$cols = @("A","B","E","F","G","H","M","Q","R")

foreach ($col in $cols) {
    $addr1 = "$col$row1"
    $addr2 = "$col$row2"
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value = $v2
    $ws.Range($addr2).Value = $v1
}
